# Update new data crawl
# This script reproduces the "update new data crawl" edit on Sheet2 of the
# TripAdvisor restaurant crawl workbook:
#  - 11 stale/duplicate review rows are removed from the data table
#    (rows 3, 4, 7, 13, 14, 15, 22, 23, 24, 25, 30 in the original 1-based
#    row numbering), shrinking the table from 315 to 304 rows.
#  - The freshly-crawled rating columns (D = nha_hang, E = an_uong, and for
#    row 2 also B = giai_tri) for the first block of rows now carry
#    non-zero scores instead of the placeholder 0.
#  - The ExternalData_1 defined name and the query table range shrink to
#    match the new row count.
#  - The sheet selection/view is moved back up near the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# 1) Remove the 11 rows whose reviews were dropped from this crawl.
#    Delete from the bottom up so earlier row numbers stay valid.
$rowsToDelete = @(30, 25, 24, 23, 22, 15, 14, 13, 7, 4, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# 2) Refresh the rating values that came back from the new crawl for the
#    first 19 review rows (rows 2-20 after the deletions above).
$ratings = @{
    2  = @{ B = 5; D = 5; E = 5 }
    3  = @{ D = 3 }
    4  = @{ D = 4; E = 5 }
    5  = @{ D = 5; E = 5 }
    6  = @{ D = 5; E = 5 }
    7  = @{ D = 5; E = 5 }
    8  = @{ D = 5; E = 5 }
    9  = @{ D = 5; E = 5 }
    10 = @{ D = 5; E = 5 }
    11 = @{ D = 5; E = 5 }
    12 = @{ D = 5; E = 5 }
    13 = @{ D = 5; E = 5 }
    14 = @{ D = 3; E = 4 }
    15 = @{ D = 1; E = 1 }
    16 = @{ D = 5; E = 5 }
    17 = @{ D = 5; E = 5 }
    18 = @{ D = 5; E = 5 }
    19 = @{ D = 5; E = 5 }
    20 = @{ D = 5 }
}

foreach ($rowNum in $ratings.Keys) {
    $cols = $ratings[$rowNum]
    foreach ($colLetter in $cols.Keys) {
        $ws.Range("$colLetter$rowNum").Value = $cols[$colLetter]
    }
}

# 3) Shrink the hidden ExternalData_1 named range to the new extent of the
#    query table (A1:G304 instead of A1:G315).
$definedName = $wb.Names.Item("Sheet2!ExternalData_1")
$definedName.RefersTo = "=Sheet2!`$A`$1:`$G`$304"

# 4) Move the view back near the top and park the selection on C17, as in
#    the saved workbook state.
$ws.Range("C17").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
